$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Poland II Liga")

# New rows of match/odds data appended to the bottom of the sheet.
$rows = @(
    @{ Row = 239; A = 237; B = 6808049; C = "Poland II Liga"; D = "Poland II Liga"; E = 45387.58333333334; F = "Skra Czestochowa";      G = "Zaglebie Lubin II";  K = 2.15; L = 3.4; M = 3;    N = 2.1;  O = 3.4; P = 3.1;  Q = -0.25; R = 1.85;  S = 1.95;  T = 2.5;  U = 2; V = 1.8; W = 0; X = 0; Y = 0; Z = 0; AA = 0 },
    @{ Row = 240; A = 238; B = 6808741; C = "Poland II Liga"; D = "Poland II Liga"; E = 45387.58333333334; F = "Olimpia Elblag";        G = "LKS Lodz II";        K = 2.3;  L = 3.2; M = 2.9;  N = 2.45; O = 3.1; P = 2.75; Q = 0;     R = 1.825; S = 1.975; T = 2.5;  U = 2; V = 1.8; W = 0; X = 0; Y = 0; Z = 0; AA = 0 },
    @{ Row = 241; A = 239; B = 6808739; C = "Poland II Liga"; D = "Poland II Liga"; E = 45387.65625;        F = "Kotwica Kolobrzeg";    G = "MKP Pogon Siedlce";  K = 2.05; L = 3.4; M = 3.2;  N = 2.05; O = 3.4; P = 3.25; Q = -0.25; R = 1.8;   S = 2;     T = 2.75; U = 2; V = 1.8; W = 0; X = 0; Y = 0; Z = 0; AA = 0 }
)

foreach ($r in $rows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.A          # A - id
    $ws.Cells.Item($row, 1).Style = "Normal"
    $ws.Range("A$row").Borders.LineStyle = 1
    $ws.Range("A$row").Font.Bold = $true
    $ws.Range("A$row").HorizontalAlignment = -4108
    $ws.Range("A$row").VerticalAlignment = -4160

    $ws.Cells.Item($row, 2).Value = $r.B          # B - match id
    $ws.Cells.Item($row, 3).Value = $r.C          # C - Div
    $ws.Cells.Item($row, 4).Value = $r.D          # D - Div Original Name

    $ws.Cells.Item($row, 5).Value = $r.E          # E - Date
    $ws.Cells.Item($row, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($row, 6).Value = $r.F          # F - HomeTeam
    $ws.Cells.Item($row, 7).Value = $r.G          # G - AwayTeam

    $ws.Cells.Item($row, 11).Value = $r.K         # K - oddH_op
    $ws.Cells.Item($row, 12).Value = $r.L         # L - oddD_op
    $ws.Cells.Item($row, 13).Value = $r.M         # M - oddA_op
    $ws.Cells.Item($row, 14).Value = $r.N         # N - oddH
    $ws.Cells.Item($row, 15).Value = $r.O         # O - oddD
    $ws.Cells.Item($row, 16).Value = $r.P         # P - oddA
    $ws.Cells.Item($row, 17).Value = $r.Q         # Q - Ah
    $ws.Cells.Item($row, 18).Value = $r.R         # R - oddAHH
    $ws.Cells.Item($row, 19).Value = $r.S         # S - oddAHA
    $ws.Cells.Item($row, 20).Value = $r.T         # T - AhOU
    $ws.Cells.Item($row, 21).Value = $r.U         # U - oddAHOver
    $ws.Cells.Item($row, 22).Value = $r.V         # V - oddAHUnder
    $ws.Cells.Item($row, 23).Value = $r.W         # W - PLH
    $ws.Cells.Item($row, 24).Value = $r.X         # X - PLD
    $ws.Cells.Item($row, 25).Value = $r.Y         # Y - PLA
    $ws.Cells.Item($row, 26).Value = $r.Z         # Z - PL_Ahh
    $ws.Cells.Item($row, 27).Value = $r.AA        # AA - PL_Aha
}
